# Update benchmark: 2026-01-05 06:47:52 UTC
# Shifts the "30,46 TL - 60,94 TL - 609,43 TL" type values one column to the
# left in several rows (F/G swap and I/K swap), updates some benchmark
# figures, and clears a couple of now-empty cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 ---
$ws.Range("F3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("G3").Value = ""
$ws.Range("I3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("K3").Value = ""

# --- Row 4 ---
$ws.Range("F4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("G4").Value = ""
$ws.Range("I4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("K4").Value = ""

# --- Row 5 ---
$ws.Range("F5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("G5").Value = ""
$ws.Range("I5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("K5").Value = ""

# --- Row 6 ---
$ws.Range("G6").Value = ""
$ws.Range("I6").Value = "6,09 TL - 12,19 TL - 152,35 TL"
$ws.Range("K6").Value = ""

# --- Row 8 ---
$ws.Range("F8").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("G8").Value = ""
$ws.Range("I8").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("K8").Value = ""

# --- Row 9 ---
$ws.Range("F9").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("G9").Value = ""
$ws.Range("I9").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("K9").Value = ""

# --- Row 10 ---
$ws.Range("F10").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("G10").Value = ""
$ws.Range("I10").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("K10").Value = ""

# --- Row 11 ---
$ws.Range("G11").Value = ""
$ws.Range("H11").Value = "3,05 TL - 6,1 TL - 76,18 TL"
$ws.Range("I11").Value = "3,04 TL - 6,09 TL - 76,17 TL"
$ws.Range("K11").Value = ""

# --- Row 12 ---
$ws.Range("G12").Value = ""
$ws.Range("K12").Value = ""

# --- Row 13 ---
$ws.Range("C13").Value = "Hesaba: Asgari 0 TL | Azami 9.999.999.999.999 TL"
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 8.700 TL"
$ws.Range("F13").Value = "Hesaba: Asgari 300 TL | Azami 3.080 TL"
$ws.Range("I13").Value = "Hesaba: Asgari 1 TL | Azami 6,09 TL"
$ws.Range("K13").Value = ""

# --- Row 14 ---
$ws.Range("F14").Value = "1.952,38 TL - 9.523,81 TL"
$ws.Range("G14").Value = ""
$ws.Range("K14").Value = ""
